# Turn the "LANGUAGE" mini-section of the CV into a "SKILLS" section:
#   - TextBox 31 (section heading): "LANGUAGE" -> "SKILLS" (and a touch shorter box)
#   - TextBox 32 (first bullet):    "Native English." -> three bullet lines:
#         "Leadership", "Critical Thinking", "Public Speech "
#   - TextBox 33 (second bullet, "Advanced spanish.") is deleted outright; the
#     two new bullet lines it would have held now live inside TextBox 32.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- TextBox 31: heading text + new (shrunk) height ---------------------
$headingShape = $s.Shapes.Item("TextBox 31")
$headingShape.TextFrame.TextRange.Text = "SKILLS"
$headingShape.Height = 179536 / 12700

# --- TextBox 32: replace single line with three bullet paragraphs -------
$bodyShape = $s.Shapes.Item("TextBox 32")
$bodyShape.TextFrame.TextRange.Text = "Leadership`rCritical Thinking`rPublic Speech "
$bodyShape.Height = 709553 / 12700

# --- TextBox 33: no longer needed, its text merged into TextBox 32 above -
$oldSpanishShape = $s.Shapes.Item("TextBox 33")
$oldSpanishShape.Delete()
